# daily auto push: 2026-01-28 22:45 UTC
# Insert two new log rows for 2026/01/29 (Thursday) just above the existing
# 2026/12/29 block, shifting all following rows down by two positions
# (old row 731 -> new row 733, ..., old row 772 -> new row 774).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at position 731; this pushes rows 731..772 down to 733..774
$ws.Rows.Item(731).Insert()
$ws.Rows.Item(731).Insert()

# Force column A to be stored as plain text so the date-like string isn't
# auto-converted into a date serial number (matches the rest of the sheet,
# which stores dates as inline text).
$ws.Range("A731:A732").NumberFormat = "@"

# Row 731: 2026/01/29, 木, 3, 201
$ws.Range("A731").Value = "2026/01/29"
$ws.Range("B731").Value = "木"
$ws.Range("C731").Value = 3
$ws.Range("D731").Value = 201

# Row 732: 2026/01/29, 木, 6, 201
$ws.Range("A732").Value = "2026/01/29"
$ws.Range("B732").Value = "木"
$ws.Range("C732").Value = 6
$ws.Range("D732").Value = 201
